# Database Window and Viewer fixed, also added Delete Record button.
# The leading "ID" column (old column A, an internal row id) is removed
# from the exported report, so every other column shifts one place to
# the left. A couple of header labels are reworded and the sample data
# row is refreshed with the new test values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old leading "ID" column - everything else shifts left by one.
$ws.Columns("A").Delete()

# --- Re-word the couple of headers that differ from a plain shift ---
$ws.Range("B1").Value = "Collection Date"
$ws.Range("T1").Value = "Recommendations"

# --- Refresh the sample data row with the new values from the report ---
# Text-valued cells: force text so numeric-looking strings (ids, phone
# numbers, dates-as-text) are not reinterpreted as numbers.
$textCells = @("A2", "B2", "E2", "G2", "I2", "J2", "T2")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("A2").Value = "2568"
$ws.Range("B2").Value = "20-03-2024"
$ws.Range("C2").Value = 265
$ws.Range("D2").Value = 265
$ws.Range("E2").Value = "asdasdasd"
$ws.Range("F2").Value = 52
$ws.Range("G2").Value = "Male"
$ws.Range("H2").Value = 23
$ws.Range("I2").Value = "32656"
$ws.Range("J2").Value = "3265623265"
$ws.Range("K2").Value = 5
$ws.Range("L2").Value = 56
$ws.Range("M2").Value = 65
$ws.Range("N2").Value = 323
$ws.Range("O2").Value = 3
$ws.Range("P2").Value = 33
$ws.Range("Q2").Value = 32
$ws.Range("R2").Value = 23
$ws.Range("S2").Value = 0.4566527413524125
$ws.Range("T2").Value = "Millets(Pearl Millet, Sorghum), Maize, Soybean, Groundnut"
